$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 4.5
$ws.Range("L7").Value = 4.75
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 8.5
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.73
$ws.Range("Y7").Value = 9
$ws.Range("AG7").Value = 11
$ws.Range("AH7").Value = 21
$ws.Range("AI7").Value = 15
$ws.Range("AM7").Value = 301
$ws.Range("AO7").Value = 10
$ws.Range("AP7").Value = 23
$ws.Range("AX7").Value = 23
$ws.Range("BA7").Value = 101
$ws.Range("BB7").Value = 251

# Row 10
$ws.Range("G10").Value = 2.55
$ws.Range("J10").Value = 3.25
$ws.Range("Y10").Value = 10
$ws.Range("AD10").Value = 6
$ws.Range("AJ10").Value = 34

# Row 11
$ws.Range("N11").Value = 9

# Row 14
$ws.Range("G14").Value = 1.8
$ws.Range("I14").Value = 4.75
$ws.Range("J14").Value = 2.4
$ws.Range("Q14").Value = 1.97
$ws.Range("R14").Value = 1.93
$ws.Range("AK14").Value = 34

# Row 16
$ws.Range("G16").Value = 1.8
$ws.Range("J16").Value = 2.6
$ws.Range("AE16").Value = 21
$ws.Range("AL16").Value = 51
$ws.Range("AQ16").Value = 41

# Row 17
$ws.Range("G17").Value = 2.05
$ws.Range("H17").Value = 3.1
$ws.Range("I17").Value = 3.8
$ws.Range("J17").Value = 2.88
$ws.Range("L17").Value = 4.5
$ws.Range("M17").Value = 1.1
$ws.Range("N17").Value = 7
$ws.Range("W17").Value = 6
$ws.Range("X17").Value = 8.5
$ws.Range("Z17").Value = 19
$ws.Range("AA17").Value = 21
$ws.Range("AB17").Value = 34
$ws.Range("AG17").Value = 8.5
$ws.Range("AH17").Value = 17
$ws.Range("AI17").Value = 13
$ws.Range("AK17").Value = 34
$ws.Range("AL17").Value = 41
$ws.Range("AM17").Value = 1250
$ws.Range("AN17").Value = 4
$ws.Range("AO17").Value = 12
$ws.Range("AU17").Value = 9
$ws.Range("AW17").Value = 5.5
$ws.Range("AX17").Value = 23
$ws.Range("AY17").Value = 34

# Row 19
$ws.Range("G19").Value = 2.25
$ws.Range("I19").Value = 3.8
$ws.Range("J19").Value = 3.2
$ws.Range("K19").Value = 1.8

